$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: new activity log entry
$ws.Range("B39").Value = 6977
$ws.Range("C39").Value = 43925
$ws.Range("D39").Value = 0.98749999999999993
$ws.Range("E39").Value = 0.010416666666666666
$ws.Range("G39").Value = "Added overview sections of LogicUnit to report. Noticed that Truth table documented is different than VHDL code"

# Row 40: new activity log entry
$ws.Range("B40").Value = 6977
$ws.Range("C40").Value = 43926
$ws.Range("D40").Value = 0.010416666666666666
$ws.Range("G40").Value = "Re-compiled VHDL code. Reproduced all diagrams to match documentation"

# Update active selection to reflect where the user ended up (A40)
$ws.Range("A40").Select()
